$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (50) down to the new rows (51-55)
$ws.Range("A50:V50").Copy()
$ws.Range("A51:V55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 51
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "azerbaijan"
$ws.Range("C51").Value = "premier-league"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45227.58333333334
$ws.Range("F51").Value = "Zira"
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = "Sabail"
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 2
$ws.Range("K51").Value = "27/10/2023 02:12"
$ws.Range("L51").Value = 1.87
$ws.Range("M51").Value = "28/10/2023 13:36"
$ws.Range("N51").Value = 2.98
$ws.Range("O51").Value = "27/10/2023 02:12"
$ws.Range("P51").Value = 3.4
$ws.Range("Q51").Value = "28/10/2023 13:36"
$ws.Range("R51").Value = 3.57
$ws.Range("S51").Value = "27/10/2023 02:12"
$ws.Range("T51").Value = 4.14
$ws.Range("U51").Value = "28/10/2023 13:36"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/zira-fk-sabail/bPZwFlyB/"

# Row 52
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "azerbaijan"
$ws.Range("C52").Value = "premier-league"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("E52").Value = 45227.70833333334
$ws.Range("F52").Value = "Turan"
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = "Kapaz"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2.03
$ws.Range("K52").Value = "27/10/2023 05:12"
$ws.Range("L52").Value = 2.06
$ws.Range("M52").Value = "28/10/2023 16:46"
$ws.Range("N52").Value = 3.1
$ws.Range("O52").Value = "27/10/2023 05:12"
$ws.Range("P52").Value = 3.4
$ws.Range("Q52").Value = "28/10/2023 16:46"
$ws.Range("R52").Value = 3.45
$ws.Range("S52").Value = "27/10/2023 05:12"
$ws.Range("T52").Value = 3.42
$ws.Range("U52").Value = "28/10/2023 16:46"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/turan-kapaz/KYzzGUL4/"

# Row 53
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "azerbaijan"
$ws.Range("C53").Value = "premier-league"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("E53").Value = 45228.54166666666
$ws.Range("F53").Value = "Sumqayit"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Gabala"
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 3.02
$ws.Range("K53").Value = "28/10/2023 02:13"
$ws.Range("L53").Value = 2.93
$ws.Range("M53").Value = "29/10/2023 12:46"
$ws.Range("N53").Value = 2.96
$ws.Range("O53").Value = "28/10/2023 02:13"
$ws.Range("P53").Value = 3.06
$ws.Range("Q53").Value = "29/10/2023 12:41"
$ws.Range("R53").Value = 2.25
$ws.Range("S53").Value = "28/10/2023 02:13"
$ws.Range("T53").Value = 2.48
$ws.Range("U53").Value = "29/10/2023 12:46"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sumqayit-fk-gabala/fiTnDS5N/"

# Row 54
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "azerbaijan"
$ws.Range("C54").Value = "premier-league"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("E54").Value = 45228.64583333334
$ws.Range("F54").Value = "Neftci Baku"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = "Sabah Baku"
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = 2.77
$ws.Range("K54").Value = "28/10/2023 04:42"
$ws.Range("L54").Value = 2.88
$ws.Range("M54").Value = "29/10/2023 14:14"
$ws.Range("N54").Value = 2.96
$ws.Range("O54").Value = "28/10/2023 04:42"
$ws.Range("P54").Value = 2.93
$ws.Range("Q54").Value = "29/10/2023 14:14"
$ws.Range("R54").Value = 2.41
$ws.Range("S54").Value = "28/10/2023 04:42"
$ws.Range("T54").Value = 2.62
$ws.Range("U54").Value = "29/10/2023 14:14"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-sabah-baku/G0UrE8jH/"

# Row 55
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "azerbaijan"
$ws.Range("C55").Value = "premier-league"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("E55").Value = 45229.66666666666
$ws.Range("F55").Value = "Araz"
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = "Qarabag"
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 5.23
$ws.Range("K55").Value = "29/10/2023 14:42"
$ws.Range("L55").Value = 5.57
$ws.Range("M55").Value = "30/10/2023 15:58"
$ws.Range("N55").Value = 3.55
$ws.Range("O55").Value = "29/10/2023 14:42"
$ws.Range("P55").Value = 3.86
$ws.Range("Q55").Value = "30/10/2023 15:58"
$ws.Range("R55").Value = 1.58
$ws.Range("S55").Value = "29/10/2023 14:42"
$ws.Range("T55").Value = 1.58
$ws.Range("U55").Value = "30/10/2023 15:58"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/araz-pfk-qarabag-agdam/QByWGA6b/"
